$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (SCD0329 -> SCD0024)
$ws.Name = "SCD0024"

# 2. Update TC_ID value in B2 (DGS-344 -> SCD0024-008)
$ws.Range("B2").Value2 = "SCD0024-008"

# 3. Re-align row 2 to the left (keeping existing vertical centering)
#    A2 (RUN): center -> left
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").HorizontalAlignment = -4131

#    B2 (TC_ID): add left + vertical-center alignment, then copy the same
#    format onto C2 (TEST_SCENARIO_DESC) so both share one style record
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

#    E2 (EXPECTED_RESULT): add left alignment (keeps vertical-center + wrap)
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").HorizontalAlignment = -4131

#    I2 (SIDEBAR_MENU): add left alignment (keeps vertical-center + wrap)
$ws.Range("I2").VerticalAlignment = -4108
$ws.Range("I2").HorizontalAlignment = -4131

#    O2: previously-empty cell now picks up the same left/center alignment
#    as its row-2 neighbours (J2, N2, P2, ...)
$ws.Range("O2").VerticalAlignment = -4108
$ws.Range("O2").HorizontalAlignment = -4131

$excel.CutCopyMode = 0

# 4. Widen column B
$ws.Columns.Item(2).ColumnWidth = 11.5

# 5. Update the active selection to B3 (and reset the scroll position)
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B3").Select()
